$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.915.62"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.544.44"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'152.72"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").Value = "3.543.49"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "4.141.62"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'32.09"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "3.542.92"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "67.668.63"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "'15.27"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "'9.71"
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("D22").Value = "'448.20"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").Value = "'0.625"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "'77.45"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("E25").Value = "  +5.63%  "
$ws.Range("D26").Value = "3.685.05"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'10.26"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'8.71"
$ws.Range("E29").Value = "  +4.66%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'1.63"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "'0.168"
$ws.Range("E32").Value = "  +7.22%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "'25.90"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "3.531.02"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'2.23"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("D42").Value = "'175.36"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "'0.0896"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("D45").Value = "'0.888"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "'28.89"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").Value = "'45.59"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'2.69"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'1.29"
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("D50").Value = "'7.62"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'0.996"
$ws.Range("E51").Value = "  -4.01%  "
